$wb = $excel.ActiveWorkbook

# --- README sheet updates ---
$wsReadme = $wb.Worksheets.Item("README")
$wsReadme.Range("A6").Value = "File Created: 2025-12-05 12:55 pm CST"
$wsReadme.Range("C8").Value = "Data Updated: 2025-12-04"

# --- Quarterly sheet updates ---
$wsQ = $wb.Worksheets.Item("Quarterly")

# Update row 92 and 93 with new data
$wsQ.Range("A92").Value = 45839
$wsQ.Range("B92").Value = 1.8
$wsQ.Range("A93").Value = 45994
$wsQ.Range("B93").Value = 1.82

# Delete old rows 94-97 (shift cells up)
$wsQ.Range("A94:B97").Delete(-4162) | Out-Null

# Update selection to reflect new end of data
$wsQ.Activate()
$wsQ.Range("B94").Select() | Out-Null
